$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Add the new row of data first (New1/New2/New3 land in the shared-string
# table before the corrected header/name strings below).
$ws1.Range("A7").Value = "New1"
$ws1.Range("B7").Value = "New2"
$ws1.Range("C7").Value = "New3"

# Fix casing / typos on existing header/data cells.
$ws1.Range("C1").Value = "Header3"
$ws1.Range("A1").Value = "Header1"
$ws1.Range("C3").Value = "Name3"

# Update selection to C3
$ws1.Range("C3").Select()
